$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data Dictionary content corrections (GMS Data Release 1):
#  - Row 3, Field column: "family_id" -> "referral_id"
#  - Row 16, Field column: "assembly" -> "genome_build"
$ws.Range("B3").Value = "referral_id"
$ws.Range("B16").Value = "genome_build"

# Restore the author's final selection/scroll position in the sheet.
$excel.Goto($ws.Range("A15"), $true)
$ws.Range("B16").Select()
